# edit.ps1 - applies the "improved version" changes to
# ling_var_multiple_ref2.xlsx:
#   - F2 default value 70 -> 75
#   - two new linguistic-variable rows appended (47: ROLE_SUBORDINATE,
#     48: SPEED_VOICE), growing the used range to A1:F48
#   - selection / scroll position moved to reflect the newly-added rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing cell tweak -------------------------------------------------
$ws.Range("F2").Value = 75

# --- new row 47: ROLE_SUBORDINATE ---------------------------------------
$ws.Range("A47").Value = "ROLE_SUBORDINATE"
$ws.Range("B47").Value = "The role of the robot is subordinate"
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = "low%%high"
$ws.Range("F47").Value = 0

# --- new row 48: SPEED_VOICE ---------------------------------------------
$ws.Range("A48").Value = "SPEED_VOICE"
$ws.Range("B48").Value = "The speed of the voice"
$ws.Range("C48").Value = 50
$ws.Range("D48").Value = 400
$ws.Range("E48").Value = "low_speed%%mid_speed%%high_speed"
$ws.Range("F48").Value = 100

# --- view state: scroll so row 22 is at the top, select the last new cell
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E48").Select()
